# ---------------------------------------------------------------------------
# Adds 2022-Q3 data:
#   1. A brand-new worksheet "2022-Q3" is inserted right before "2022-Q2"
#      (which pushes 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3 / 2021-Q2 one
#      position later in the tab order; their own content is untouched).
#   2. The "总计" (summary) sheet gets a new row for 2022-Q3 inserted right
#      after the header row, with the older rows shifting down and their
#      running index (column A) renumbered.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: create the "2022-Q3" worksheet ahead of "2022-Q2".
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($q2Sheet)
$ws.Name = "2022-Q3"

# Match the look & feel (header / index-column styling) of the existing
# quarterly sheets by pulling the formatting (not the values) from 2022-Q2.
# NB: re-fetch the "2022-Q2" handle now that the sheet collection has
# changed shape (Add() above can invalidate previously-captured refs).
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Range("A1:H2").Copy()
$ws.Range("A1:H2").PasteSpecial(-4122)

# The fund-code / ratio-looking columns must stay text (e.g. "012186" must
# not collapse to 12186, "0.3730" must not collapse to 0.373), so force
# those columns to a text format before writing the values.
$ws.Range("B2:G18").NumberFormat = "@"

$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "012186"
$ws.Cells.Item(2,3).Value = "招商品质成长混合A"
$ws.Cells.Item(2,4).Value = "12.27"
$ws.Cells.Item(2,5).Value = "91.08"
$ws.Cells.Item(2,6).Value = "3.04"
$ws.Cells.Item(2,7).Value = "0.3730"
$ws.Cells.Item(2,8).Value = 10
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "011855"
$ws.Cells.Item(3,3).Value = "银华长荣混合"
$ws.Cells.Item(3,4).Value = "10.55"
$ws.Cells.Item(3,5).Value = "65.92"
$ws.Cells.Item(3,6).Value = "3.01"
$ws.Cells.Item(3,7).Value = "0.3176"
$ws.Cells.Item(3,8).Value = 7
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "630010"
$ws.Cells.Item(4,3).Value = "华商价值精选混合"
$ws.Cells.Item(4,4).Value = "4.30"
$ws.Cells.Item(4,5).Value = "81.81"
$ws.Cells.Item(4,6).Value = "5.56"
$ws.Cells.Item(4,7).Value = "0.2391"
$ws.Cells.Item(4,8).Value = 1
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "003291"
$ws.Cells.Item(5,3).Value = "信澳健康中国灵活配置混合A"
$ws.Cells.Item(5,4).Value = "7.25"
$ws.Cells.Item(5,5).Value = "91.02"
$ws.Cells.Item(5,6).Value = "2.87"
$ws.Cells.Item(5,7).Value = "0.2081"
$ws.Cells.Item(5,8).Value = 10
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "009360"
$ws.Cells.Item(6,3).Value = "招商创新增长混合A"
$ws.Cells.Item(6,4).Value = "5.73"
$ws.Cells.Item(6,5).Value = "91.71"
$ws.Cells.Item(6,6).Value = "3.21"
$ws.Cells.Item(6,7).Value = "0.1839"
$ws.Cells.Item(6,8).Value = 9
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "011598"
$ws.Cells.Item(7,3).Value = "信澳医药健康混合"
$ws.Cells.Item(7,4).Value = "5.45"
$ws.Cells.Item(7,5).Value = "88.25"
$ws.Cells.Item(7,6).Value = "3.30"
$ws.Cells.Item(7,7).Value = "0.1798"
$ws.Cells.Item(7,8).Value = 8
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "012187"
$ws.Cells.Item(8,3).Value = "招商品质成长混合C"
$ws.Cells.Item(8,4).Value = "5.25"
$ws.Cells.Item(8,5).Value = "91.08"
$ws.Cells.Item(8,6).Value = "3.04"
$ws.Cells.Item(8,7).Value = "0.1596"
$ws.Cells.Item(8,8).Value = 10
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "015208"
$ws.Cells.Item(9,3).Value = "信澳健康中国灵活配置混合C"
$ws.Cells.Item(9,4).Value = "3.72"
$ws.Cells.Item(9,5).Value = "91.02"
$ws.Cells.Item(9,6).Value = "2.87"
$ws.Cells.Item(9,7).Value = "0.1068"
$ws.Cells.Item(9,8).Value = 10
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "008978"
$ws.Cells.Item(10,3).Value = "银华长丰混合"
$ws.Cells.Item(10,4).Value = "2.24"
$ws.Cells.Item(10,5).Value = "73.86"
$ws.Cells.Item(10,6).Value = "3.92"
$ws.Cells.Item(10,7).Value = "0.0878"
$ws.Cells.Item(10,8).Value = 6
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "630006"
$ws.Cells.Item(11,3).Value = "华商产业升级混合"
$ws.Cells.Item(11,4).Value = "0.85"
$ws.Cells.Item(11,5).Value = "81.97"
$ws.Cells.Item(11,6).Value = "5.61"
$ws.Cells.Item(11,7).Value = "0.0477"
$ws.Cells.Item(11,8).Value = 1
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "009361"
$ws.Cells.Item(12,3).Value = "招商创新增长混合C"
$ws.Cells.Item(12,4).Value = "0.85"
$ws.Cells.Item(12,5).Value = "91.71"
$ws.Cells.Item(12,6).Value = "3.21"
$ws.Cells.Item(12,7).Value = "0.0273"
$ws.Cells.Item(12,8).Value = 9
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "015032"
$ws.Cells.Item(13,3).Value = "中融医药消费混合A"
$ws.Cells.Item(13,4).Value = "0.54"
$ws.Cells.Item(13,5).Value = "90.81"
$ws.Cells.Item(13,6).Value = "3.99"
$ws.Cells.Item(13,7).Value = "0.0215"
$ws.Cells.Item(13,8).Value = 7
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "005493"
$ws.Cells.Item(14,3).Value = "鑫元价值精选灵活配置混合A"
$ws.Cells.Item(14,4).Value = "0.55"
$ws.Cells.Item(14,5).Value = "76.82"
$ws.Cells.Item(14,6).Value = "2.94"
$ws.Cells.Item(14,7).Value = "0.0162"
$ws.Cells.Item(14,8).Value = 9
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "006193"
$ws.Cells.Item(15,3).Value = "鑫元核心资产股票A"
$ws.Cells.Item(15,4).Value = "0.11"
$ws.Cells.Item(15,5).Value = "83.48"
$ws.Cells.Item(15,6).Value = "3.05"
$ws.Cells.Item(15,7).Value = "0.0034"
$ws.Cells.Item(15,8).Value = 8
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "015033"
$ws.Cells.Item(16,3).Value = "中融医药消费混合C"
$ws.Cells.Item(16,4).Value = "0.02"
$ws.Cells.Item(16,5).Value = "90.81"
$ws.Cells.Item(16,6).Value = "3.99"
$ws.Cells.Item(16,7).Value = "0.0008"
$ws.Cells.Item(16,8).Value = 7
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "005494"
$ws.Cells.Item(17,3).Value = "鑫元价值精选灵活配置混合C"
$ws.Cells.Item(17,4).Value = "0.01"
$ws.Cells.Item(17,5).Value = "76.82"
$ws.Cells.Item(17,6).Value = "2.94"
$ws.Cells.Item(17,7).Value = "0.0003"
$ws.Cells.Item(17,8).Value = 9
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "006194"
$ws.Cells.Item(18,3).Value = "鑫元核心资产股票C"
$ws.Cells.Item(18,4).Value = "0.01"
$ws.Cells.Item(18,5).Value = "83.48"
$ws.Cells.Item(18,6).Value = "3.05"
$ws.Cells.Item(18,7).Value = "0.0003"
$ws.Cells.Item(18,8).Value = 8

# ---------------------------------------------------------------------------
# Part 2: insert the 2022-Q3 row into the "总计" summary sheet.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing quarterly rows (currently rows 2-6) down by one row.
$total.Rows.Item(2).Insert()

# The insert above bleeds the header row's formatting onto the new row;
# re-copy the plain data-row styling (from what is now row 3) so row 2
# matches the other quarter rows instead of looking like the header.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 17
$total.Cells.Item(2,4).Value = 1.97

# Renumber the running index in column A for the rows that shifted down.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
